# Generate Report for Handback
#
# The localization-status report workbook gets a new handback recorded for
# the "5475b91d-53d6-4bcb-bf91-c84cce0a7ee7" source file: row 7 of both the
# "zh-cn" and "de-de" target-language worksheets is filled in with the
# target handback file, its datetime, and an error detail explaining that
# the handed-back file isn't built from the latest source revision.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# zh-cn sheet, row 7 (source row for 5475b91d-53d6-4bcb-bf91-c84cce0a7ee7)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

# "Latest Target File" (column I) becomes a hyperlink to the handback .md,
# mirroring the pattern already used for rows 2-5 (I2..I5).
$ws2.Hyperlinks.Add($ws2.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/edc7cffae329f0624e35aee3af2bfa600db77f69/e2e/5475b91d-53d6-4bcb-bf91-c84cce0a7ee7.md", "", "", "5475b91d-53d6-4bcb-bf91-c84cce0a7ee7.md")

# "Latest Handback File" (column J)
$ws2.Range("J7").Value = "5475b91d-53d6-4bcb-bf91-c84cce0a7ee7.aebbe8d2f0d0f48524b723d4bbc0bba54ec7c683.zh-cn.xlf"

# "Latest Handback DateTime" (column K)
$ws2.Range("K7").Value = "2016-09-03 13:00:26"

# "Error Detail" (column P)
$ws2.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/edc7cffae329f0624e35aee3af2bfa600db77f69/e2e/5475b91d-53d6-4bcb-bf91-c84cce0a7ee7.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/054bad01315876da02232dc164893ae73afabf8d/e2e/5475b91d-53d6-4bcb-bf91-c84cce0a7ee7.md."

# ---------------------------------------------------------------------
# de-de sheet, row 7 (same source row, other target language)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Add($ws3.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/edc7cffae329f0624e35aee3af2bfa600db77f69/e2e/5475b91d-53d6-4bcb-bf91-c84cce0a7ee7.md", "", "", "5475b91d-53d6-4bcb-bf91-c84cce0a7ee7.md")

$ws3.Range("J7").Value = "5475b91d-53d6-4bcb-bf91-c84cce0a7ee7.aebbe8d2f0d0f48524b723d4bbc0bba54ec7c683.de-de.xlf"

$ws3.Range("K7").Value = "2016-09-03 13:00:33"

$ws3.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/edc7cffae329f0624e35aee3af2bfa600db77f69/e2e/5475b91d-53d6-4bcb-bf91-c84cce0a7ee7.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/054bad01315876da02232dc164893ae73afabf8d/e2e/5475b91d-53d6-4bcb-bf91-c84cce0a7ee7.md."
